$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 2503.375
$ws.Range("I4").Value = 1754.5
$ws.Range("K4").Value = 1754.5
$ws.Range("M4").Value = -1640.5
$ws.Range("H12").Value = 25703.125
$ws.Range("I12").Value = 1040
$ws.Range("J12").Value = 66808.336
$ws.Range("K12").Value = 1040
$ws.Range("L12").Value = 66808.336
$ws.Range("M12").Value = -870
$ws.Range("N12").Value = -67148.336
$ws.Range("H70").Value = 14290520
$ws.Range("I70").Value = 40005496
$ws.Range("K70").Value = 120016488
$ws.Range("M70").Value = -120016218
$ws.Range("H73").Value = 14290520
$ws.Range("I73").Value = 40005496
$ws.Range("K73").Value = 120016488
$ws.Range("M73").Value = -120015552
$ws.Range("H112").Value = 78860.53999999999
$ws.Range("J112").Value = 102258.7
$ws.Range("L112").Value = 306776.1
$ws.Range("N112").Value = -308992.1
$ws.Range("H119").Value = 2000
$ws.Range("J119").Value = 2000
$ws.Range("L119").Value = 6000
$ws.Range("N119").Value = -15676
$ws.Range("H132").Value = 2151.1875
$ws.Range("I132").Value = 2027.4918
$ws.Range("K132").Value = 6082.4754
$ws.Range("M132").Value = -3552.4754
$ws.Range("H137").Value = 9023.031000000001
$ws.Range("I137").Value = 11890.909
$ws.Range("J137").Value = 2713.7
$ws.Range("K137").Value = 35672.727
$ws.Range("L137").Value = 8141.099999999999
$ws.Range("M137").Value = -33122.727
$ws.Range("N137").Value = -13241.1
$ws.Range("H138").Value = 3076.4092
$ws.Range("J138").Value = 5310.3447
$ws.Range("L138").Value = 15931.0341
$ws.Range("N138").Value = -26211.0341

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4192.0576
$ws.Range("I32").Value = 3846.6736
$ws.Range("K32").Value = 3846.6736
$ws.Range("M32").Value = -3559.6736
$ws.Range("H92").Value = 220034700
$ws.Range("J92").Value = 275032000
$ws.Range("L92").Value = 275032000
$ws.Range("N92").Value = -275036992
$ws.Range("H96").Value = 41335.25
$ws.Range("J96").Value = 41335.25
$ws.Range("L96").Value = 41335.25
$ws.Range("N96").Value = -46827.25
$ws.Range("H97").Value = 26277.857
$ws.Range("I97").Value = 14234.375
$ws.Range("J97").Value = 42335.832
$ws.Range("K97").Value = 14234.375
$ws.Range("L97").Value = 42335.832
$ws.Range("M97").Value = -13738.375
$ws.Range("N97").Value = -43327.832

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4024.5527
$ws.Range("I31").Value = 3480.2856
$ws.Range("J31").Value = 5548.5
$ws.Range("K31").Value = 3480.2856
$ws.Range("L31").Value = 5548.5
$ws.Range("M31").Value = -3185.2856
$ws.Range("N31").Value = -6138.5
$ws.Range("H34").Value = 4024.5527
$ws.Range("I34").Value = 3480.2856
$ws.Range("J34").Value = 5548.5
$ws.Range("K34").Value = 3480.2856
$ws.Range("L34").Value = 5548.5
$ws.Range("M34").Value = -3278.2856
$ws.Range("N34").Value = -5952.5
$ws.Range("H86").Value = 15999.75
$ws.Range("H89").Value = 15999.75
$ws.Range("H99").Value = 16592873
$ws.Range("I99").Value = 29031528
$ws.Range("K99").Value = 29031528
$ws.Range("M99").Value = -29030030
$ws.Range("H107").Value = 30310506
$ws.Range("I107").Value = 47630144
$ws.Range("J107").Value = 1142.5834
$ws.Range("K107").Value = 47630144
$ws.Range("L107").Value = 1142.5834
$ws.Range("M107").Value = -47628224
$ws.Range("N107").Value = -4982.5834
$ws.Range("H126").Value = 16592873
$ws.Range("I126").Value = 29031528
$ws.Range("K126").Value = 87094584
$ws.Range("M126").Value = -87092114

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 256
$ws.Range("I14").Value = 256
$ws.Range("K14").Value = 768
$ws.Range("M14").Value = -595
$ws.Range("H64").Value = 2337
$ws.Range("I64").Value = 1725
$ws.Range("K64").Value = 5175
$ws.Range("M64").Value = -4905
$ws.Range("H67").Value = 2337
$ws.Range("I67").Value = 1725
$ws.Range("K67").Value = 5175
$ws.Range("M67").Value = -4239
$ws.Range("H107").Value = 682.61536
$ws.Range("I107").Value = 273.8
$ws.Range("J107").Value = 742.7353000000001
$ws.Range("K107").Value = 821.4000000000001
$ws.Range("L107").Value = 2228.2059
$ws.Range("M107").Value = 1098.6
$ws.Range("N107").Value = -6068.2059

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 33333.168
$ws.Range("J57").Value = 39999.5
$ws.Range("L57").Value = 39999.5
$ws.Range("N57").Value = -41639.5
$ws.Range("H59").Value = 9494.5
$ws.Range("I59").Value = 8000
$ws.Range("J59").Value = 10989
$ws.Range("K59").Value = 8000
$ws.Range("L59").Value = 10989
$ws.Range("M59").Value = -7417
$ws.Range("N59").Value = -12155
$ws.Range("H70").Value = 10603.36
$ws.Range("I70").Value = 9566
$ws.Range("J70").Value = 13270.857
$ws.Range("K70").Value = 9566
$ws.Range("L70").Value = 13270.857
$ws.Range("M70").Value = -9296
$ws.Range("N70").Value = -13810.857
$ws.Range("H73").Value = 10603.36
$ws.Range("I73").Value = 9566
$ws.Range("J73").Value = 13270.857
$ws.Range("K73").Value = 9566
$ws.Range("L73").Value = 13270.857
$ws.Range("M73").Value = -8630
$ws.Range("N73").Value = -15142.857

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 22596.28
$ws.Range("I7").Value = 41446.184
$ws.Range("K7").Value = 41446.184
$ws.Range("M7").Value = -41334.184
$ws.Range("H22").Value = 7153.622
$ws.Range("I22").Value = 8370.166999999999
$ws.Range("K22").Value = 8370.166999999999
$ws.Range("M22").Value = -8075.166999999999
$ws.Range("H27").Value = 7153.622
$ws.Range("I27").Value = 8370.166999999999
$ws.Range("K27").Value = 8370.166999999999
$ws.Range("M27").Value = -8263.166999999999
$ws.Range("H40").Value = 29171.85
$ws.Range("I40").Value = 40150.273
$ws.Range("K40").Value = 40150.273
$ws.Range("M40").Value = -40014.273
$ws.Range("H46").Value = 2166.875
$ws.Range("I46").Value = 1547.9166
$ws.Range("K46").Value = 1547.9166
$ws.Range("M46").Value = -1359.9166
$ws.Range("H122").Value = 6884.3823
$ws.Range("I122").Value = 5657.5864
$ws.Range("K122").Value = 16972.7592
$ws.Range("M122").Value = -14522.7592
$ws.Range("H126").Value = 22596.28
$ws.Range("I126").Value = 41446.184
$ws.Range("K126").Value = 124338.552
$ws.Range("M126").Value = -121868.552
$ws.Range("H132").Value = 375721.06
$ws.Range("I132").Value = 498446.78
$ws.Range("K132").Value = 1495340.34
$ws.Range("M132").Value = -1492810.34
$ws.Range("H136").Value = 8810.65
$ws.Range("I136").Value = 2370.6667
$ws.Range("J136").Value = 14079.728
$ws.Range("K136").Value = 7112.000100000001
$ws.Range("L136").Value = 42239.18399999999
$ws.Range("M136").Value = -4562.000100000001
$ws.Range("N136").Value = -47339.18399999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 9509.5
$ws.Range("I33").Value = 8019
$ws.Range("J33").Value = 11000
$ws.Range("K33").Value = 8019
$ws.Range("L33").Value = 11000
$ws.Range("M33").Value = -7769
$ws.Range("N33").Value = -11500
$ws.Range("H36").Value = 9509.5
$ws.Range("I36").Value = 8019
$ws.Range("J36").Value = 11000
$ws.Range("K36").Value = 8019
$ws.Range("L36").Value = 11000
$ws.Range("M36").Value = -7769
$ws.Range("N36").Value = -11500
$ws.Range("H96").Value = 20002400
$ws.Range("I96").Value = 33335334
$ws.Range("K96").Value = 33335334
$ws.Range("M96").Value = -33333961
$ws.Range("H136").Value = 2358.7144
$ws.Range("I136").Value = 1480.3914
$ws.Range("J136").Value = 6399
$ws.Range("K136").Value = 4441.174199999999
$ws.Range("L136").Value = 19197
$ws.Range("M136").Value = -1891.174199999999
$ws.Range("N136").Value = -24297
